$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("D Green")
$ws.Range("B2").Value = 8.264462809917356
$ws.Range("B3").Value = 5.952380952380952
$ws.Range("B10").Value = 6.018518518518518
$ws.Range("B14").Value = 7.725321888412018

$ws = $wb.Worksheets.Item("Green")
$ws.Range("B2").Value = 13.22314049586777
$ws.Range("B3").Value = 27.77777777777778
$ws.Range("B4").Value = 23.04347826086957
$ws.Range("B5").Value = 12.5
$ws.Range("B6").Value = 15.94827586206897
$ws.Range("B7").Value = 55.81395348837209
$ws.Range("B8").Value = 19.04761904761905
$ws.Range("B9").Value = 60.8
$ws.Range("B10").Value = 27.77777777777778
$ws.Range("B11").Value = 19.81981981981982
$ws.Range("B12").Value = 19.04761904761905
$ws.Range("B13").Value = 26.78571428571428
$ws.Range("B14").Value = 21.88841201716738
$ws.Range("B15").Value = 38.75968992248063

$ws = $wb.Worksheets.Item("Yellow")
$ws.Range("B2").Value = 35.9504132231405
$ws.Range("B3").Value = 19.84126984126984
$ws.Range("B4").Value = 47.39130434782609
$ws.Range("B5").Value = 43.53448275862069
$ws.Range("B6").Value = 62.5
$ws.Range("B7").Value = 20.34883720930233
$ws.Range("B8").Value = 25.59523809523809
$ws.Range("B9").Value = 17.2
$ws.Range("B10").Value = 31.48148148148148
$ws.Range("B11").Value = 40.54054054054054
$ws.Range("B12").Value = 25.59523809523809
$ws.Range("B13").Value = 42.85714285714285
$ws.Range("B14").Value = 42.91845493562232
$ws.Range("B15").Value = 27.51937984496124

$ws = $wb.Worksheets.Item("Orange")
$ws.Range("B2").Value = 19.00826446280992
$ws.Range("B3").Value = 12.3015873015873
$ws.Range("B4").Value = 7.82608695652174
$ws.Range("B5").Value = 18.96551724137931
$ws.Range("B6").Value = 8.620689655172415
$ws.Range("B7").Value = 11.62790697674419
$ws.Range("B8").Value = 30.95238095238095
$ws.Range("B9").Value = 6.4
$ws.Range("B10").Value = 15.27777777777778
$ws.Range("B11").Value = 15.76576576576577
$ws.Range("B12").Value = 30.95238095238095
$ws.Range("B13").Value = 8.928571428571429
$ws.Range("B14").Value = 11.58798283261803
$ws.Range("B15").Value = 12.01550387596899

$ws = $wb.Worksheets.Item("Brown")
$ws.Range("B2").Value = 6.198347107438017
$ws.Range("B3").Value = 13.49206349206349
$ws.Range("B4").Value = 9.565217391304348
$ws.Range("B5").Value = 8.620689655172415
$ws.Range("B6").Value = 4.310344827586207
$ws.Range("B7").Value = 4.651162790697674
$ws.Range("B8").Value = 10.11904761904762
$ws.Range("B9").Value = 4
$ws.Range("B10").Value = 1.851851851851852
$ws.Range("B11").Value = 7.207207207207207
$ws.Range("B12").Value = 10.11904761904762
$ws.Range("B13").Value = 9.821428571428571
$ws.Range("B14").Value = 6.437768240343347
$ws.Range("B15").Value = 5.426356589147287

$ws = $wb.Worksheets.Item("Red")
$ws.Range("B2").Value = 8.264462809917356
$ws.Range("B3").Value = 16.66666666666667
$ws.Range("B4").Value = 6.086956521739131
$ws.Range("B5").Value = 11.20689655172414
$ws.Range("B6").Value = 3.448275862068965
$ws.Range("B7").Value = 5.232558139534884
$ws.Range("B8").Value = 10.71428571428571
$ws.Range("B9").Value = 7.199999999999999
$ws.Range("B10").Value = 8.333333333333332
$ws.Range("B11").Value = 12.16216216216216
$ws.Range("B12").Value = 10.71428571428571
$ws.Range("B13").Value = 11.60714285714286
$ws.Range("B14").Value = 5.150214592274678
$ws.Range("B15").Value = 12.4031007751938

$ws = $wb.Worksheets.Item("Default Red")
$ws.Range("B2").Value = 5.785123966942149
$ws.Range("B3").Value = 3.968253968253968
$ws.Range("B4").Value = 6.086956521739131
$ws.Range("B5").Value = 5.172413793103448
$ws.Range("B6").Value = 5.172413793103448
$ws.Range("B7").Value = 2.325581395348837
$ws.Range("B8").Value = 3.571428571428571
$ws.Range("B9").Value = 3.2
$ws.Range("B10").Value = 9.25925925925926
$ws.Range("B11").Value = 4.504504504504505
$ws.Range("B12").Value = 3.571428571428571
$ws.Range("B14").Value = 4.291845493562231
$ws.Range("B15").Value = 3.875968992248062

$ws = $wb.Worksheets.Item("Blue")
$ws.Range("B2").Value = 3.305785123966942
$ws.Range("B3").Value = 0
$ws.Range("B4").Value = 0
$ws.Range("B5").Value = 0
$ws.Range("B6").Value = 0
$ws.Range("B7").Value = 0
$ws.Range("B8").Value = 0
$ws.Range("B9").Value = 1.2
$ws.Range("B10").Value = 0
$ws.Range("B11").Value = 0
$ws.Range("B12").Value = 0
$ws.Range("B13").Value = 0
$ws.Range("B14").Value = 0
$ws.Range("B15").Value = 0
